# Update the Ascend_Roadmap sheet: shift the roadmap table from a
# time-series-by-row layout to a time-series-by-column layout, and
# clear out the now-unused last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ascend_Roadmap")

# Row 3: column headers -> quarter/year labels.
# C3:F3 look like plain numbers ("2025".."2028"); force them to be stored
# as text (matching the rest of the table) instead of being auto-coerced
# to numeric values, then drop the temporary formatting so no stray style
# is left behind on the cells.
$ws.Range("A3").Value = "2024-Q1"
$ws.Range("B3").Value = "2024-Q4"
$ws.Range("C3:F3").NumberFormat = "@"
$ws.Range("C3").Value = "2025"
$ws.Range("D3").Value = "2026"
$ws.Range("E3").Value = "2027"
$ws.Range("F3").Value = "2028"
$ws.Range("C3:F3").ClearFormats()

# Row 4: product names
$ws.Range("A4").Value = "Ascend 910B"
$ws.Range("B4").Value = "Ascend 910C"
$ws.Range("C4").Value = "Ascend 950"
$ws.Range("D4").Value = "Ascend 960"
$ws.Range("E4").Value = "Ascend 970"
$ws.Range("F4").Value = "下一代"

# Row 5: process node
$ws.Range("A5").Value = "7nm (国内)"
$ws.Range("B5").Value = "7nm (国内)"
$ws.Range("C5").Value = "5nm (国内)"
$ws.Range("D5").Value = "5nm (国内)"
$ws.Range("E5").Value = "3nm (国内)"
$ws.Range("F5").Value = "2nm (国内)"

# Row 6: AI performance
$ws.Range("A6").Value = "FP16 400 TFLOPS"
$ws.Range("B6").Value = "FP16 ~500 TFLOPS"
$ws.Range("C6").Value = "FP16 ~1 PFLOPS"
$ws.Range("D6").Value = "FP16 ~2 PFLOPS"
$ws.Range("E6").Value = "FP16 ~4 PFLOPS"
$ws.Range("F6").Value = "FP4 ~8 ZettaFLOPS"

# Row 7: key features
$ws.Range("A7").Value = "当前主力"
$ws.Range("B7").Value = "性能提升"
$ws.Range("C7").Value = "新一代"
$ws.Range("D7").Value = "下一代"
$ws.Range("E7").Value = "下一代"
$ws.Range("F7").Value = "远期目标"

# Row 8: remarks/status
$ws.Range("A8").Value = "已发布"
$ws.Range("B8").Value = "已发布"
$ws.Range("C8").Value = "2026发布"
$ws.Range("D8").Value = "2027目标"
$ws.Range("E8").Value = "2028目标"
$ws.Range("F8").Value = "2028目标"

# Row 9 is no longer used - clear any remaining content from the old layout,
# but keep the (now empty) row present in the sheet.
$ws.Range("A9:F9").ClearContents()
$ws.Rows.Item(9).OutlineLevel = 0
